$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header title strings (Volume Number + date range) ---
$ws.Range("A8").Value = "Volume 30   Number  40"
$ws.Range("C9").Value = "Report Covering the Week  10/2/2023  Through  10/8/2023"

# --- Update crime-statistics grid (rows 14-30) ---
$ws.Range("L14").Value = -100
$ws.Range("L14").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("D15").Value = 1
$ws.Range("D15").NumberFormat = "#,##0"
$ws.Range("E15").Value = -100
$ws.Range("E15").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("G15").Value = 1
$ws.Range("G15").NumberFormat = "#,##0"
$ws.Range("H15").Value = 0
$ws.Range("H15").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("J15").Value = 17
$ws.Range("K15").Value = -23.529411764705
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = -85.714285714285
$ws.Range("F16").Value = 5
$ws.Range("G16").Value = 16
$ws.Range("H16").Value = -68.75
$ws.Range("I16").Value = 109
$ws.Range("J16").Value = 150
$ws.Range("K16").Value = -27.333333333333
$ws.Range("L16").Value = 15.957446808510
$ws.Range("M16").Value = -46.039603960396
$ws.Range("N16").Value = -84.797768479776
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = -50
$ws.Range("F17").Value = 22
$ws.Range("G17").Value = 19
$ws.Range("H17").Value = 15.789473684210
$ws.Range("I17").Value = 177
$ws.Range("J17").Value = 149
$ws.Range("K17").Value = 18.791946308724
$ws.Range("L17").Value = 51.282051282051
$ws.Range("M17").Value = 132.894736842105
$ws.Range("N17").Value = -7.329842931937
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 9
$ws.Range("E18").Value = -66.666666666666
$ws.Range("F18").Value = 9
$ws.Range("G18").Value = 29
$ws.Range("H18").Value = -68.965517241379
$ws.Range("I18").Value = 205
$ws.Range("J18").Value = 179
$ws.Range("K18").Value = 14.525139664804
$ws.Range("L18").Value = 40.410958904109
$ws.Range("M18").Value = -5.963302752293
$ws.Range("N18").Value = -81.874447391688
$ws.Range("C19").Value = 11
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = 22.222222222222
$ws.Range("F19").Value = 47
$ws.Range("G19").Value = 52
$ws.Range("H19").Value = -9.615384615384
$ws.Range("I19").Value = 471
$ws.Range("J19").Value = 515
$ws.Range("K19").Value = -8.543689320388
$ws.Range("L19").Value = 72.527472527472
$ws.Range("M19").Value = 42.296072507552
$ws.Range("N19").Value = 0.212765957446
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 25
$ws.Range("G20").Value = 27
$ws.Range("H20").Value = -7.407407407407
$ws.Range("I20").Value = 280
$ws.Range("J20").Value = 199
$ws.Range("K20").Value = 40.703517587939
$ws.Range("L20").Value = 127.642276422764
$ws.Range("M20").Value = 56.424581005586
$ws.Range("N20").Value = -92.625757176718
$ws.Range("C21").Value = 22
$ws.Range("D21").Value = 35
$ws.Range("E21").Value = -37.142857142857
$ws.Range("F21").Value = 109
$ws.Range("G21").Value = 144
$ws.Range("H21").Value = -24.305555555555
$ws.Range("I21").Value = 1255
$ws.Range("J21").Value = 1209
$ws.Range("K21").Value = 3.804797353184
$ws.Range("L21").Value = 62.564766839378
$ws.Range("M21").Value = 23.280943025540
$ws.Range("N21").Value = -80.195676187470
$ws.Range("C22").Value = 1
$ws.Range("C22").NumberFormat = "#,##0"
$ws.Range("E22").Value = -50
$ws.Range("F22").Value = 3
$ws.Range("G22").Value = 6
$ws.Range("I22").Value = 34
$ws.Range("J22").Value = 25
$ws.Range("K22").Value = 36
$ws.Range("L22").Value = 6.25
$ws.Range("M22").Value = 112.5
$ws.Range("C23").Value = 1
$ws.Range("C23").NumberFormat = "#,##0"
$ws.Range("I23").Value = 61
$ws.Range("K23").Value = 52.5
$ws.Range("L23").Value = 96.774193548387
$ws.Range("M23").Value = 134.615384615385
$ws.Range("C24").Value = 31
$ws.Range("D24").Value = 27
$ws.Range("E24").Value = 14.814814814814
$ws.Range("F24").Value = 113
$ws.Range("G24").Value = 115
$ws.Range("H24").Value = -1.739130434782
$ws.Range("I24").Value = 1115
$ws.Range("J24").Value = 1113
$ws.Range("K24").Value = 0.179694519317
$ws.Range("L24").Value = 69.452887537993
$ws.Range("M24").Value = 43.685567010309
$ws.Range("C25").Value = 16
$ws.Range("E25").Value = 128.571428571429
$ws.Range("F25").Value = 44
$ws.Range("H25").Value = 46.666666666666
$ws.Range("I25").Value = 369
$ws.Range("J25").Value = 346
$ws.Range("K25").Value = 6.647398843930
$ws.Range("L25").Value = 42.471042471042
$ws.Range("M25").Value = 19.032258064516
$ws.Range("D26").Value = 1
$ws.Range("D26").NumberFormat = "#,##0"
$ws.Range("E26").Value = -100
$ws.Range("E26").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = -50
$ws.Range("J26").Value = 26
$ws.Range("K26").Value = -15.384615384615
$ws.Range("D27").Value = 1
$ws.Range("D27").NumberFormat = "#,##0"
$ws.Range("E27").Value = 0
$ws.Range("E27").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F27").Value = 3
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = -50
$ws.Range("I27").Value = 35
$ws.Range("J27").Value = 49
$ws.Range("K27").Value = -28.571428571428
$ws.Range("L27").Value = -12.5
$ws.Range("L28").Value = -60
$ws.Range("L29").Value = -60
$ws.Range("D30").Value = 1
$ws.Range("D30").NumberFormat = "#,##0"
$ws.Range("E30").Value = -100
$ws.Range("E30").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("G30").Value = 1
$ws.Range("G30").NumberFormat = "#,##0"
$ws.Range("H30").Value = 100
$ws.Range("H30").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("J30").Value = 4
$ws.Range("K30").Value = 425

# --- C30: text "0" sharing style with D30 (numeric style -> text style) ---
$ws.Range("D30").Copy() | Out-Null
$ws.Range("C30").PasteSpecial(-4104) | Out-Null
$ws.Range("D30").Copy() | Out-Null
$ws.Range("C30").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
